$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Deploy to AWS -> set Priority (column B) to 5
$ws.Range("B7").Value = 5

# Row 8: Scale Heatmaps task is now complete
# Set Effort (column B) to 4
$ws.Range("B8").Value = 4
# Set Status (column D) to "Complete"
$ws.Range("D8").Value = "Complete"
# Apply the gray fill style used by other completed/prioritized rows (copy from row 6 which has same style pattern)
$ws.Range("A6:E6").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122) # xlPasteFormats

# Update selection to B7 as shown in the diff
$ws.Range("B7").Select()
